# V0R2: Estilizacion del banner y sobre
#
# Applies spell-checker proofErr markers (splitting runs around words
# Word's proofing tool flagged) and adds a few new list items that were
# introduced in this revision (a Camel_case wiki link under "Anexos:" and
# several new CSS/vocabulary entries under "Vocabulario:").

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-ParagraphXml($searchText, $xml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Paragraph not found for search text: $searchText"
    }
    $para = $rng.Paragraphs(1)
    [void]$para.Range.InsertXML($xml)
}

function Insert-ParagraphAfter($searchText, $xml) {
    # Locates the paragraph containing $searchText and inserts a brand new
    # sibling paragraph immediately after it (built from $xml), by inserting
    # just before that paragraph's own end-of-paragraph mark.
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor paragraph not found for search text: $searchText"
    }
    $para = $rng.Paragraphs(1)
    $insertPoint = $d.Range($para.Range.End - 1, $para.Range.End - 1)
    [void]$insertPoint.InsertXML($xml)
}

# 1. Title: "Arquitectura CSS: Descomplicando los problemas"
Replace-ParagraphXml "Arquitectura CSS: Descomplicando los problemas" (
    "<w:p $wNs><w:pPr><w:pStyle w:val=`"Ttulo1`"/></w:pPr>" +
    "<w:r><w:t xml:space=`"preserve`">Arquitectura CSS: </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>Descomplicando</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> los problemas</w:t></w:r>" +
    "</w:p>"
)

# 2. "Conclusión 1: Layout base y estilización de cabecera"
Replace-ParagraphXml "Conclusión 1: Layout base y estilización de cabecera" (
    "<w:p $wNs><w:pPr><w:pStyle w:val=`"Ttulo2`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr>" +
    "<w:r><w:t>Conclusión 1:</w:t></w:r>" +
    "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>Layout</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> base y estilización de cabecera</w:t></w:r>" +
    "</w:p>"
)

# 3. "Organizar y estructurar los archivos css del proyecto;"
Replace-ParagraphXml "Organizar y estructurar los archivos css del proyecto;" (
    "<w:p $wNs><w:pPr><w:pStyle w:val=`"Prrafodelista`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr></w:pPr>" +
    "<w:r><w:t xml:space=`"preserve`">Organizar y estructurar los archivos </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>css</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> del proyecto;</w:t></w:r>" +
    "</w:p>"
)

# 4. "Organización y estructura de archivos .css;"
Replace-ParagraphXml "Organización y estructura de archivos .css;" (
    "<w:p $wNs><w:pPr><w:pStyle w:val=`"Prrafodelista`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"4`"/></w:numPr></w:pPr>" +
    "<w:r><w:t>Organización y estructura de archivos .</w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>css</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t>;</w:t></w:r>" +
    "</w:p>"
)

# 5. "Las imágenes ilustrativas deben estar referenciadas en nuestros archivos .css;"
Replace-ParagraphXml "Las imágenes ilustrativas deben estar referenciadas en nuestros archivos .css;" (
    "<w:p $wNs><w:pPr><w:pStyle w:val=`"Prrafodelista`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"5`"/></w:numPr></w:pPr>" +
    "<w:r><w:t>Las imágenes ilustrativas deben estar referenciadas en nuestros archivos .</w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>css</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t>;</w:t></w:r>" +
    "</w:p>"
)

# 6. "Conclusión 5: Estilización del pie de página responsividad"
Replace-ParagraphXml "Conclusión 5: Estilización del pie de página responsividad" (
    "<w:p $wNs><w:pPr><w:pStyle w:val=`"Ttulo2`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr>" +
    "<w:r><w:t>Conclusión 5:</w:t></w:r>" +
    "<w:r><w:t xml:space=`"preserve`"> Estilización del pie de página </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>responsividad</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "</w:p>"
)

# 7. New paragraph: Camel_case wikipedia link, inserted right after the
#    "bradfrost.com" Anexos link (before the "Código:" heading).
Insert-ParagraphAfter "https://bradfrost.com/blog/post/atomic-web-design/" (
    "<w:p $wNs><w:pPr><w:pStyle w:val=`"Prrafodelista`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"7`"/></w:numPr></w:pPr>" +
    "<w:r><w:t>https://es.wikipedia.org/wiki/Camel_case</w:t></w:r>" +
    "</w:p>"
)

# 8. ".menu-link:hover{  text-decoration: underline;"
Replace-ParagraphXml ".menu-link:hover{  text-decoration: underline;" (
    "<w:p $wNs><w:pPr><w:pStyle w:val=`"Prrafodelista`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"7`"/></w:numPr></w:pPr>" +
    "<w:r><w:t>.</w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>menu-link:hover</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`">{  </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>text-decoration</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`">: </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>underline</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t>;</w:t></w:r>" +
    "</w:p>"
)

# 9. "text-transform: lowercase;"
Replace-ParagraphXml "text-transform: lowercase;" (
    "<w:p $wNs><w:pPr><w:pStyle w:val=`"Prrafodelista`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"7`"/></w:numPr></w:pPr>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>text-transform</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`">: </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>lowercase</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t>;</w:t></w:r>" +
    "</w:p>"
)

# 10. "Atomic desing"
Replace-ParagraphXml "Atomic desing" (
    "<w:p $wNs><w:pPr><w:pStyle w:val=`"Prrafodelista`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"7`"/></w:numPr></w:pPr>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>Atomic</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>desing</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "</w:p>"
)

# 11-15. New vocabulary list items appended after "Atomic desing", in order.
Insert-ParagraphAfter "Atomic desing" (
    "<w:p $wNs><w:pPr><w:pStyle w:val=`"Prrafodelista`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"7`"/></w:numPr></w:pPr>" +
    "<w:r><w:t>Camel case</w:t></w:r>" +
    "</w:p>"
)

Insert-ParagraphAfter "Camel case" (
    "<w:p $wNs><w:pPr><w:pStyle w:val=`"Prrafodelista`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"7`"/></w:numPr></w:pPr>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>Bem</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`">: Block </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>element</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>modifier</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "</w:p>"
)

Insert-ParagraphAfter "Bem: Block element modifier" (
    "<w:p $wNs><w:pPr><w:pStyle w:val=`"Prrafodelista`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"7`"/></w:numPr></w:pPr>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>height</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`">: </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>calc</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t>(100vh - 72px);</w:t></w:r>" +
    "</w:p>"
)

Insert-ParagraphAfter "height: calc(100vh - 72px);" (
    "<w:p $wNs><w:pPr><w:pStyle w:val=`"Prrafodelista`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"7`"/></w:numPr></w:pPr>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>text-shadow</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`">: 0 4px </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>4px</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>rgba</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t>(0,0,0,0.75)</w:t></w:r>" +
    "</w:p>"
)

Insert-ParagraphAfter "text-shadow: 0 4px 4px rgba(0,0,0,0.75)" (
    "<w:p $wNs><w:pPr><w:pStyle w:val=`"Prrafodelista`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"7`"/></w:numPr></w:pPr>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>transform</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`">: </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>translate</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t>(-50%, -50%);</w:t></w:r>" +
    "</w:p>"
)

Write-Host "Edits applied."
